$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TESTDATA")

# Add new row 16: LETTER_A / a
$ws.Range("A16").Value = "LETTER_A"
$ws.Range("B16").Value = "a"

# Change B3 from "Alice" to "incorrect"
$ws.Range("B3").Value = "incorrect"

# Update the conditional formatting range from A11:A15 to A11:A16 to include the new row
$cf = $ws.Range("A11:A15").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("A11:A16"))

# Update the selection to B17
$ws.Range("B17").Select()
